$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.139.85'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.607.37'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.06%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '203.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +10.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '566.59'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.66%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.604.51'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.623'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.28%  '
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("E10").Value = '  +0.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '61.18'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +14.85%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.152'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.37%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000288'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +10.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.185.19'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.610.53'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.25%  '
$ws.Range("E17").Value = '  +0.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.00'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '67.987.21'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.43%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.39'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.97%  '
$ws.Range("E21").Value = '  +1.84%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '402.86'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.25'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +15.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.43'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.81%  '
$ws.Range("E26").Value = '  +0.90%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.64'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.12%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.94'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +10.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.35'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +17.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.42'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.79%  '
$ws.Range("E32").Value = '  +1.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '670.26'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '12.25'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.40%  '
$ws.Range("E35").Value = '  +0.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '63.74'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '42.21'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.72%  '
$ws.Range("E38").Value = '  +4.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.314.97'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0₃0770'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.31%  '
$ws.Range("E42").Value = '  +12.21%  '
$ws.Range("E43").Value = '  +3.82%  '
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.78'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +9.42%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.05'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +30.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.997'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0420'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.76'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +11.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.87'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.132'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.40%  '
$ws.Range("E51").Value = '  -0.31%  '
